$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GS")
$ws.Activate()

$ws.Range("F2").Value = 1003
$ws.Range("F3").Value = 1003
$ws.Range("F4").Value = 1003

$ws.Range("H2").Value = "JWMT"
$ws.Range("H3").Value = "JWMT"
$ws.Range("H4").Value = "JWMT"

$ws.Range("F5").Select()
